# Add a new "Croatia" market sheet, cloned from the existing "Turkey" sheet
# (same layout/formatting), then fill in the Croatia-specific values and
# hand the "active sheet / selected tab" state over to the new sheet - this
# mirrors the way every other per-market sheet in this workbook was created.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Clone the Turkey sheet (keeps column widths, merged cells, styles, etc.)
# and drop the copy immediately after it, i.e. at the end of the tab strip.
$turkey.Copy([System.Reflection.Missing]::Value, $turkey)
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in the market-specific values (these become new shared-string entries).
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2415/T2416/T2418"

# Column B needs to be a bit wider to comfortably fit the longer Jira key.
$croatia.Columns("B:B").ColumnWidth = 31.3

# Restore Turkey's selection back to a plain "whole sheet" selection and
# hand the active-tab / tab-selected flag over to the freshly added sheet,
# exactly like the previous "last sheet" used to hold it.
$turkey.Activate()
$turkey.Cells.Select() | Out-Null

$croatia.Activate()
$croatia.Range("A1:D1").Select() | Out-Null
